# Add a "Różnica godzin" (hours difference) column (G):
#   Różnica godzin = Suma godzin (D) - "Idealne" godziny (F)
# including closed days (where both D and F are already 0, giving a 0 difference).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell G2: new label, formatted like the other header cells (F2) ---
$ws.Range("G2").Value = "Różnica godzin"
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # xlPasteFormats - copy F2's style onto G2

# --- Data rows 4-33: G = D - F (difference between actual and "ideal" hours) ---
for ($r = 4; $r -le 33; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 7).Value = $d - $f
}

# Give the new data cells (G4:G33) the same look as the neighbouring F column
$ws.Range("F4:F33").Copy()
$ws.Range("G4:G33").PasteSpecial(-4122)   # xlPasteFormats

# --- Totals row 34: overall difference is (by design) zero ---
$ws.Range("G34").Value = 0.0
$ws.Range("F34").Copy()
$ws.Range("G34").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Column width for the new column ---
$ws.Columns.Item(7).ColumnWidth = 13.7

# --- Number formats: tighten decimal precision on the two numeric formats in use ---
$ws.Range("D4:D33").NumberFormat = "#.#"
$ws.Range("F4:F33").NumberFormat = "#.#"
$ws.Range("G4:G33").NumberFormat = "#.#"

$ws.Range("A34:F34").NumberFormat = "#"
$ws.Range("G34").NumberFormat = "#"
